# The handoff transform failed for the source file in both the zh-cn and
# de-de target sheets, so the generated .xlf handoff file/link, its
# datetime, and the "Include" reason no longer apply - the row reverts to
# the "nothing has happened yet" state and the status explains why.
$wb = $excel.ActiveWorkbook

$sheetNames = @("zh-cn", "de-de")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    # Remember every hyperlink except the one on the "Latest Handoff File"
    # cell (C2) - that link goes away because there is no handoff file
    # anymore.
    $keep = @()
    foreach ($hl in $ws.Hyperlinks) {
        $addr = $hl.Range.Address()
        if ($addr -ne '$C$2') {
            $keep += , @($addr, $hl.Address, $hl.TextToDisplay)
        }
    }

    # Status: "Ready for handoff" -> "Handoff transform failed".
    $ws.Range("B2").Value = "Handoff transform failed"

    # Latest Handoff File: drop the link + the generated .xlf file name.
    $ws.Hyperlinks.Delete()
    $ws.Range("C2").Clear()

    # Latest Handoff Datetime reverts to the unset sentinel datetime.
    $ws.Range("D2").Value = "0001-01-01 00:00:00"

    # Latest Handback DateTime also reports the unset sentinel datetime.
    $ws.Range("G2").Value = "0001-01-01 00:00:00"

    # Handoff Reason: "Include" -> "Ignored".
    $ws.Range("H2").Value = "Ignored"

    # Restore the hyperlinks (and their visual style) that should remain.
    foreach ($item in $keep) {
        $ws.Hyperlinks.Add($ws.Range($item[0]), $item[1], "", "", $item[2])
        $ws.Range($item[0]).Style = "HyperLink"
    }
}
